$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the 12 months of data (rows 2-13) twice more, into rows 14-25 and 26-37
$src = $ws.Range("A2:G13")
$values = $src.Value2

$dest1 = $ws.Range("A14:G25")
$dest1.Value2 = $values

$dest2 = $ws.Range("A26:G37")
$dest2.Value2 = $values

# Match the existing date style on column A (style used by A2:A13) without
# introducing a new number format entry
$ws.Range("A2").Copy()
$ws.Range("A14:A25").PasteSpecial(-4122)
$ws.Range("A26:A37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("I23").Select()
